$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.646.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.59%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.857.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.51%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.034"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +3.01%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'322.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.74%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'1.030"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.63%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4403"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.57%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3795"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.73%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.07427"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +2.69%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.8802"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.88%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'21.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.20%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.871.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -8.07%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'5.539"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.85%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'6.720"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.33%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.07220"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.76%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'83.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.31%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'1.036"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.16%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.000009080"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.78%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'1.030"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.69%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'15.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.44%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'27.662.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.51%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'5.288"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.85%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'11.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +4.28%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'158.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.93%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'1.923"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.43%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'18.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.61%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'1.984"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +4.03%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'5.305"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.35%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'117.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.32%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.09070"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.53%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.208"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +4.04%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.7647"
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").Value = "'4.547"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.67%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'2.888"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +3.27%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'1.031"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.27%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'1.156"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.69%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.01980"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +3.02%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.05329"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.22%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.5181"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.87%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'2.830"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.57%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.1686"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.32%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'6.778"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +5.44%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'8.634"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +4.15%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'109.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.06%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'10.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.31%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'1.720"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +4.15%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.4669"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.85%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.06417"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.23%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'1.859"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.91%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'39.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +4.76%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'64.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.92%  "
$ws.Range("E51").Style = "Normal"
